$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: was a formula "=0" with style s=2; becomes literal number 5.5, default style (s=0)
$ws.Range("C3").ClearFormats()
$ws.Range("C3").Value = 5.5

# D3: new boolean TRUE cell with a custom "TRUE/FALSE" number format (numFmtId 165)
$ws.Range("D3").Value = $true
$ws.Range("D3").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Rows 7-9 (B,C,D): re-apply formatting (style moves from s=0 to s=1)
$ws.Range("B7:D9").Style = "Normal"

# Move the active selection from C3 to E3
$ws.Range("E3").Select()
